$wb = $excel.ActiveWorkbook
$ws = $wb.Sheets.Item(1)

# Rename the sheet to reflect the new "through" date
$ws.Name = "Through 2022-06-26"

# Update the June row label to the new "through" date
$ws.Range("A7").Value = "June (through 06-26)"

# Update June row (row 7) values
$ws.Range("B7").Value = 15
$ws.Range("C7").Value = 31
$ws.Range("D7").Value = 64
$ws.Range("E7").Value = 49
$ws.Range("F7").Value = 36
$ws.Range("G7").Value = 98
$ws.Range("H7").Value = 104
$ws.Range("I7").Value = 125

# Update Total row (row 8) values
$ws.Range("B8").Value = 123
$ws.Range("C8").Value = 240
$ws.Range("D8").Value = 380
$ws.Range("E8").Value = 344
$ws.Range("F8").Value = 240
$ws.Range("G8").Value = 456
$ws.Range("H8").Value = 735
$ws.Range("I8").Value = 788
